$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handback for e2e\62956314-... file failed to transform: status text changes everywhere
# it is used - Overview row 3 (E3/F3) as well as the Status column (C3) on each language sheet.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Error Detail column (P) widened to fit the new long messages (to match the
# existing 40-wide columns elsewhere on the sheet, e.g. column A).
$wsZhCn.Columns.Item(16).ColumnWidth = $wsZhCn.Columns.Item(1).ColumnWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $wsDeDe.Columns.Item(1).ColumnWidth

# Populate the new Error Detail messages for the failed handback row (row 3).
$wsZhCn.Range("P3").Value = "Handback file name: u31mgn0y.1t3 is different with handoff file name: 62956314-3861-4101-b7d3-8a412f437fa1.33265c254e41853231e21966b85ce62e769411bd.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: u31mgn0y.1t3 is different with handoff file name: 62956314-3861-4101-b7d3-8a412f437fa1.33265c254e41853231e21966b85ce62e769411bd.de-de."
